$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 23.15092491234345
$ws.Range("C2").Value = 11.50559536387624
$ws.Range("D2").Value = 3.794202529467003
$ws.Range("E2").Value = 9.520947567715242
$ws.Range("F2").Value = 55.22734134379357
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 39.01110950967163
$ws.Range("J2").Value = 9.626041788798206
$ws.Range("L2").Value = 12.08883899679955
$ws.Range("M2").Value = 20.04156972374729
$ws.Range("B3").Value = 22.95528901744763
$ws.Range("C3").Value = 11.18283497767065
$ws.Range("D3").Value = 3.750218158558777
$ws.Range("E3").Value = 9.509481565454852
$ws.Range("F3").Value = 55.09400305996552
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 38.9799580843771
$ws.Range("J3").Value = 9.632950141144429
$ws.Range("L3").Value = 12.10803170345066
$ws.Range("M3").Value = 20.03019113898947
$ws.Range("B4").Value = 22.84119477541338
$ws.Range("C4").Value = 10.98399651579488
$ws.Range("D4").Value = 3.722577827539784
$ws.Range("E4").Value = 9.502287124511396
$ws.Range("F4").Value = 55.02314259088069
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 38.96752281961452
$ws.Range("J4").Value = 9.637402168283172
$ws.Range("L4").Value = 12.12149539988678
$ws.Range("M4").Value = 20.02737034863743
$ws.Range("B5").Value = 22.79626530066256
$ws.Range("C5").Value = 10.90294522310036
$ws.Range("D5").Value = 3.711157926610346
$ws.Range("E5").Value = 9.499316276164311
$ws.Range("F5").Value = 54.99704266630791
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 38.96413604877585
$ws.Range("J5").Value = 9.639269427259176
$ws.Range("L5").Value = 12.12740444659023
$ws.Range("M5").Value = 20.02727080368228
$ws.Range("B6").Value = 22.78890068905012
$ws.Range("C6").Value = 10.88948964561648
$ws.Range("D6").Value = 3.709252315909858
$ws.Range("E6").Value = 9.498820602516366
$ws.Range("F6").Value = 54.99287665089449
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 38.96367511073037
$ws.Range("J6").Value = 9.639582690907275
$ws.Range("L6").Value = 12.12841116151964
$ws.Range("M6").Value = 20.02731774869769
$ws.Range("B7").Value = 22.84058244201054
$ws.Range("C7").Value = 10.98290332851927
$ws.Range("D7").Value = 3.722424443200116
$ws.Range("E7").Value = 9.502247217073609
$ws.Range("F7").Value = 55.02277934899333
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 38.96747034289723
$ws.Range("J7").Value = 9.63742713591158
$ws.Range("L7").Value = 12.12157338050921
$ws.Range("M7").Value = 20.0273647521922
$ws.Range("B8").Value = 23.08224927401326
$ws.Range("C8").Value = 11.39452909974629
$ws.Range("D8").Value = 3.77917010917876
$ws.Range("E8").Value = 9.517025948582843
$ws.Range("F8").Value = 55.17908666122197
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 38.99897865775903
$ws.Range("J8").Value = 9.628380236185247
$ws.Range("L8").Value = 12.09510823779843
$ws.Range("M8").Value = 20.0367834340781
$ws.Range("B9").Value = 23.60166931332714
$ws.Range("C9").Value = 12.19052311084924
$ws.Range("D9").Value = 3.88528679116595
$ws.Range("E9").Value = 9.544795592379595
$ws.Range("F9").Value = 55.57251317239209
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 39.113919785215
$ws.Range("J9").Value = 9.612300863427494
$ws.Range("L9").Value = 12.05652508193661
$ws.Range("M9").Value = 20.0881759483127
$ws.Range("B10").Value = 24.00776550211862
$ws.Range("C10").Value = 12.76112083059044
$ws.Range("D10").Value = 3.959932998069328
$ws.Range("E10").Value = 9.564478105054368
$ws.Range("F10").Value = 55.91379024267841
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 39.23080006627891
$ws.Range("J10").Value = 9.601490567596414
$ws.Range("L10").Value = 12.03628165082906
$ws.Range("M10").Value = 20.14579597212312
$ws.Range("B11").Value = 24.197089749426
$ws.Range("C11").Value = 13.01622907358608
$ws.Range("D11").Value = 3.993136144497622
$ws.Range("E11").Value = 9.573279939766051
$ws.Range("F11").Value = 56.08017240740767
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 39.29100081743056
$ws.Range("J11").Value = 9.596788444410015
$ws.Range("L11").Value = 12.02882844033613
$ws.Range("M11").Value = 20.17626483313726
$ws.Range("B12").Value = 24.26937839115243
$ws.Range("C12").Value = 13.11208809076407
$ws.Range("D12").Value = 4.005598000437521
$ws.Range("E12").Value = 9.576591426373872
$ws.Range("F12").Value = 56.14475551287068
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 39.31480525104384
$ws.Range("J12").Value = 9.595038711626696
$ws.Range("L12").Value = 12.02625817527583
$ws.Range("M12").Value = 20.18840892065321
$ws.Range("B13").Value = 24.25378427610168
$ws.Range("C13").Value = 13.09147781853349
$ws.Range("D13").Value = 4.002919121565211
$ws.Range("E13").Value = 9.575879194047051
$ws.Range("F13").Value = 56.13077660802407
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 39.30963378653603
$ws.Range("J13").Value = 9.595414177547074
$ws.Range("L13").Value = 12.02680052155135
$ws.Range("M13").Value = 20.18576661539479
$ws.Range("B14").Value = 24.20302541182297
$ws.Range("C14").Value = 13.02413094271541
$ws.Range("D14").Value = 3.994163639819009
$ws.Range("E14").Value = 9.573552808727669
$ws.Range("F14").Value = 56.08545421861058
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 39.29293904938036
$ws.Range("J14").Value = 9.59664387514283
$ws.Range("L14").Value = 12.02861193272631
$ws.Range("M14").Value = 20.17725182201433
$ws.Range("B15").Value = 24.17200979640076
$ws.Range("C15").Value = 12.98277910740366
$ws.Range("D15").Value = 3.988786050405951
$ws.Range("E15").Value = 9.572125028191643
$ws.Range("F15").Value = 56.05789769130686
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 39.28284415740836
$ws.Range("J15").Value = 9.5974011152342
$ws.Range("L15").Value = 12.02975429450697
$ws.Range("M15").Value = 20.17211502107256
$ws.Range("B16").Value = 23.99547958540363
$ws.Range("C16").Value = 12.7443500937433
$ws.Range("D16").Value = 3.957747662508102
$ws.Range("E16").Value = 9.56389985105025
$ws.Range("F16").Value = 55.9031388849148
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 39.22700702818606
$ws.Range("J16").Value = 9.601802188003425
$ws.Range("L16").Value = 12.03680403421918
$ws.Range("M16").Value = 20.14388992814098
$ws.Range("B17").Value = 23.88831391870139
$ws.Range("C17").Value = 12.59686267898442
$ws.Range("D17").Value = 3.938511347709395
$ws.Range("E17").Value = 9.558815264698691
$ws.Range("F17").Value = 55.81103520902285
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 39.19455156808121
$ws.Range("J17").Value = 9.604557211272649
$ws.Range("L17").Value = 12.04157824250376
$ws.Range("M17").Value = 20.12766105319313
$ws.Range("B18").Value = 23.82711090997115
$ws.Range("C18").Value = 12.51161766420608
$ws.Range("D18").Value = 3.927376355883014
$ws.Range("E18").Value = 9.555876437947717
$ws.Range("F18").Value = 55.75910909734066
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 39.1765461350855
$ws.Range("J18").Value = 9.606162122912108
$ws.Range("L18").Value = 12.04448951282469
$ws.Range("M18").Value = 20.11872774095454
$ws.Range("B19").Value = 23.80646542549495
$ws.Range("C19").Value = 12.48268734834907
$ws.Range("D19").Value = 3.923594178510952
$ws.Range("E19").Value = 9.55487893660813
$ws.Range("F19").Value = 55.74170872475418
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 39.17056360546833
$ws.Range("J19").Value = 9.606709008231482
$ws.Range("L19").Value = 12.04550361456421
$ws.Range("M19").Value = 20.11577213814371
$ws.Range("B20").Value = 23.89967724139782
$ws.Range("C20").Value = 12.61260661029006
$ws.Range("D20").Value = 3.940566433382981
$ws.Range("E20").Value = 9.559358006468651
$ws.Range("F20").Value = 55.82073134494212
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 39.19793800475956
$ws.Range("J20").Value = 9.604261834851206
$ws.Range("L20").Value = 12.041052917279
$ws.Range("M20").Value = 20.12934717304811
$ws.Range("B21").Value = 24.21791886977616
$ws.Range("C21").Value = 13.04393334970538
$ws.Range("D21").Value = 3.996738385848682
$ws.Range("E21").Value = 9.574236708040058
$ws.Range("F21").Value = 56.09872387913176
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 39.29781538123501
$ws.Range("J21").Value = 9.596281846707354
$ws.Range("L21").Value = 12.02807303836185
$ws.Range("M21").Value = 20.17973642210976
$ws.Range("B22").Value = 24.42934997005899
$ws.Range("C22").Value = 13.32144672920978
$ws.Range("D22").Value = 4.03279900477269
$ws.Range("E22").Value = 9.583835150586339
$ws.Range("F22").Value = 56.28959218729899
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 39.36896247380623
$ws.Range("J22").Value = 9.591246272537898
$ws.Range("L22").Value = 12.02105919641518
$ws.Range("M22").Value = 20.21619885441327
$ws.Range("B23").Value = 24.31621165351953
$ws.Range("C23").Value = 13.17376554721121
$ws.Range("D23").Value = 4.013613327755354
$ws.Range("E23").Value = 9.578723693156501
$ws.Range("F23").Value = 56.18689011804467
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 39.33045407238122
$ws.Range("J23").Value = 9.593917444844811
$ws.Range("L23").Value = 12.02466830684019
$ws.Range("M23").Value = 20.19641729361957
$ws.Range("B24").Value = 23.89453860763723
$ws.Range("C24").Value = 12.60549018052601
$ws.Range("D24").Value = 3.939637564948931
$ws.Range("E24").Value = 9.55911268161889
$ws.Range("F24").Value = 55.81634452590284
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 39.19640496054024
$ws.Range("J24").Value = 9.604395309073279
$ws.Range("L24").Value = 12.04128989813176
$ws.Range("M24").Value = 20.12858364160923
$ws.Range("B25").Value = 23.45660906970396
$ws.Range("C25").Value = 11.97718198357478
$ws.Range("D25").Value = 3.857149873532544
$ws.Range("E25").Value = 9.537410925461575
$ws.Range("F25").Value = 55.45684425001159
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 39.07712773733518
$ws.Range("J25").Value = 9.616473885024014
$ws.Range("L25").Value = 12.02680052155135
$ws.Range("M25").Value = 20.07076742212642
